$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.814.86"
$ws.Range("E2").Value = "  +5.95%  "
$ws.Range("D3").Value = "3.651.97"
$ws.Range("E3").Value = "  +6.04%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'595.08"
$ws.Range("E5").Value = "  +2.75%  "
$ws.Range("D6").Value = "'194.05"
$ws.Range("E6").Value = "  +3.57%  "
$ws.Range("D7").Value = "'0.647"
$ws.Range("E7").Value = "  +2.67%  "
$ws.Range("D8").Value = "3.644.32"
$ws.Range("E8").Value = "  +5.98%  "
$ws.Range("E10").Value = "  +7.34%  "
$ws.Range("D11").Value = "'0.673"
$ws.Range("E11").Value = "  +4.49%  "
$ws.Range("E12").Value = "  +2.86%  "
$ws.Range("D13").Value = "'0.0000294"
$ws.Range("E13").Value = "  +6.25%  "
$ws.Range("D14").Value = "'9.93"
$ws.Range("E14").Value = "  +5.73%  "
$ws.Range("D15").Value = "4.230.81"
$ws.Range("E15").Value = "  +5.95%  "
$ws.Range("D16").Value = "'20.05"
$ws.Range("E16").Value = "  +7.19%  "
$ws.Range("D17").Value = "3.652.44"
$ws.Range("E17").Value = "  +5.92%  "
$ws.Range("D18").Value = "70.751.31"
$ws.Range("E18").Value = "  +5.75%  "
$ws.Range("E19").Value = "  +6.09%  "
$ws.Range("E20").Value = "  +2.89%  "
$ws.Range("E21").Value = "  +4.49%  "
$ws.Range("D22").Value = "'489.63"
$ws.Range("E22").Value = "  +1.55%  "
$ws.Range("D23").Value = "'19.16"
$ws.Range("E23").Value = "  +13.62%  "
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("D25").Value = "'4.52"
$ws.Range("E25").Value = "  +4.46%  "
$ws.Range("D26").Value = "'91.57"
$ws.Range("E26").Value = "  +2.47%  "
$ws.Range("E27").Value = "  +6.84%  "
$ws.Range("D28").Value = "'11.48"
$ws.Range("E28").Value = "  +4.62%  "
$ws.Range("D29").Value = "'9.60"
$ws.Range("E29").Value = "  +6.22%  "
$ws.Range("D30").Value = "'32.87"
$ws.Range("E30").Value = "  +5.28%  "
$ws.Range("D31").Value = "'7.76"
$ws.Range("E31").Value = "  +6.46%  "
$ws.Range("E32").Value = "  +9.81%  "
$ws.Range("D33").Value = "'628.69"
$ws.Range("E33").Value = "  +5.29%  "
$ws.Range("D34").Value = "'12.29"
$ws.Range("E34").Value = "  +4.55%  "
$ws.Range("D35").Value = "'66.20"
$ws.Range("E35").Value = "  +3.99%  "
$ws.Range("D36").Value = "'40.23"
$ws.Range("E36").Value = "  +9.79%  "
$ws.Range("D37").Value = "'0.413"
$ws.Range("E37").Value = "  +6.79%  "
$ws.Range("D38").Value = "0.0₃0825"
$ws.Range("E38").Value = "  +9.72%  "
$ws.Range("D39").Value = "'0.149"
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("D41").Value = "'3.59"
$ws.Range("E41").Value = "  +1.45%  "
$ws.Range("D42").Value = "3.306.14"
$ws.Range("E42").Value = "  +1.80%  "
$ws.Range("D43").Value = "'3.16"
$ws.Range("E43").Value = "  +9.35%  "
$ws.Range("D44").Value = "'2.81"
$ws.Range("E44").Value = "  +11.13%  "
$ws.Range("D45").Value = "'0.0458"
$ws.Range("E45").Value = "  +6.68%  "
$ws.Range("D46").Value = "'3.04"
$ws.Range("E46").Value = "  +6.88%  "
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").Value = "'9.41"
$ws.Range("E47").Value = "  +8.27%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'3.31"
$ws.Range("E48").Value = "  +2.42%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "'0.139"
$ws.Range("E49").Value = "  +3.60%  "
$ws.Range("D50").Value = "'3.30"
$ws.Range("E50").Value = "  -1.91%  "
$ws.Range("E51").Value = "  +0.05%  "
